$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to populate: column A (IP-like numbers), column B (123), column C (whitelist name)
# Row 2 changes from "MOHITO" to "KOLA"; rows 3-13 are new.
$data = @(
    @(123,   123, "KOLA"),    # row 2 (existing row, value changes)
    @(2345,  123, "KOLA"),    # row 3
    @(124,   123, "KOLA"),    # row 4
    @(623,   123, "KOLA"),    # row 5
    @(45732, 123, "KOLA"),    # row 6
    @(5243,  123, "KOLA"),    # row 7
    @(236,   123, "KOLA"),    # row 8
    @(2365,  123, "tvhome"),  # row 9
    @(5687,  123, "KOLA"),    # row 10
    @(55,    123, "KOLA"),    # row 11
    @(456,   123, "tvhome"),  # row 12
    @(6768,  123, "KOLA")     # row 13
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("G13").Select()
